$d = $word.ActiveDocument
$tbl = $d.Tables.Item(1)

$cell = $tbl.Cell(1, 1)
$p = $cell.Range.Paragraphs.Item(1).Range
$p.InsertXML('<?xml version="1.0" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:rPr><w:sz w:val="32"/></w:rPr><w:t>26 x 67</w:t><w:br/><w:t xml:space="preserve">  6    7</w:t><w:br/><w:t xml:space="preserve">  ----</w:t><w:br/><w:t>2|    |</w:t><w:br/><w:t>6|    |</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>')

$cell = $tbl.Cell(1, 2)
$p = $cell.Range.Paragraphs.Item(1).Range
$p.InsertXML('<?xml version="1.0" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:rPr><w:sz w:val="32"/></w:rPr><w:t>14 x 72</w:t><w:br/><w:t xml:space="preserve">  7    2</w:t><w:br/><w:t xml:space="preserve">  ----</w:t><w:br/><w:t>1|    |</w:t><w:br/><w:t>4|    |</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>')

$cell = $tbl.Cell(1, 3)
$p = $cell.Range.Paragraphs.Item(1).Range
$p.InsertXML('<?xml version="1.0" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:rPr><w:sz w:val="32"/></w:rPr><w:t>21 x 67</w:t><w:br/><w:t xml:space="preserve">  6    7</w:t><w:br/><w:t xml:space="preserve">  ----</w:t><w:br/><w:t>2|    |</w:t><w:br/><w:t>1|    |</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>')

$cell = $tbl.Cell(2, 1)
$p = $cell.Range.Paragraphs.Item(1).Range
$p.InsertXML('<?xml version="1.0" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:rPr><w:sz w:val="32"/></w:rPr><w:t>40 x 97</w:t><w:br/><w:t xml:space="preserve">  9    7</w:t><w:br/><w:t xml:space="preserve">  ----</w:t><w:br/><w:t>4|    |</w:t><w:br/><w:t>0|    |</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>')

$cell = $tbl.Cell(2, 2)
$p = $cell.Range.Paragraphs.Item(1).Range
$p.InsertXML('<?xml version="1.0" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:rPr><w:sz w:val="32"/></w:rPr><w:t>32 x 57</w:t><w:br/><w:t xml:space="preserve">  5    7</w:t><w:br/><w:t xml:space="preserve">  ----</w:t><w:br/><w:t>3|    |</w:t><w:br/><w:t>2|    |</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>')

$cell = $tbl.Cell(2, 3)
$p = $cell.Range.Paragraphs.Item(1).Range
$p.InsertXML('<?xml version="1.0" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:rPr><w:sz w:val="32"/></w:rPr><w:t>90 x 71</w:t><w:br/><w:t xml:space="preserve">  7    1</w:t><w:br/><w:t xml:space="preserve">  ----</w:t><w:br/><w:t>9|    |</w:t><w:br/><w:t>0|    |</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>')

$cell = $tbl.Cell(3, 1)
$p = $cell.Range.Paragraphs.Item(1).Range
$p.InsertXML('<?xml version="1.0" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:rPr><w:sz w:val="32"/></w:rPr><w:t>12 x 27</w:t><w:br/><w:t xml:space="preserve">  2    7</w:t><w:br/><w:t xml:space="preserve">  ----</w:t><w:br/><w:t>1|    |</w:t><w:br/><w:t>2|    |</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>')

$cell = $tbl.Cell(3, 2)
$p = $cell.Range.Paragraphs.Item(1).Range
$p.InsertXML('<?xml version="1.0" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:rPr><w:sz w:val="32"/></w:rPr><w:t>48 x 50</w:t><w:br/><w:t xml:space="preserve">  5    0</w:t><w:br/><w:t xml:space="preserve">  ----</w:t><w:br/><w:t>4|    |</w:t><w:br/><w:t>8|    |</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>')

$cell = $tbl.Cell(3, 3)
$p = $cell.Range.Paragraphs.Item(1).Range
$p.InsertXML('<?xml version="1.0" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:rPr><w:sz w:val="32"/></w:rPr><w:t>50 x 43</w:t><w:br/><w:t xml:space="preserve">  4    3</w:t><w:br/><w:t xml:space="preserve">  ----</w:t><w:br/><w:t>5|    |</w:t><w:br/><w:t>0|    |</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>')

$cell = $tbl.Cell(4, 1)
$p = $cell.Range.Paragraphs.Item(1).Range
$p.InsertXML('<?xml version="1.0" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:rPr><w:sz w:val="32"/></w:rPr><w:t>76 x 90</w:t><w:br/><w:t xml:space="preserve">  9    0</w:t><w:br/><w:t xml:space="preserve">  ----</w:t><w:br/><w:t>7|    |</w:t><w:br/><w:t>6|    |</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>')

$cell = $tbl.Cell(4, 2)
$p = $cell.Range.Paragraphs.Item(1).Range
$p.InsertXML('<?xml version="1.0" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:rPr><w:sz w:val="32"/></w:rPr><w:t>88 x 94</w:t><w:br/><w:t xml:space="preserve">  9    4</w:t><w:br/><w:t xml:space="preserve">  ----</w:t><w:br/><w:t>8|    |</w:t><w:br/><w:t>8|    |</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>')

$cell = $tbl.Cell(4, 3)
$p = $cell.Range.Paragraphs.Item(1).Range
$p.InsertXML('<?xml version="1.0" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:rPr><w:sz w:val="32"/></w:rPr><w:t>83 x 35</w:t><w:br/><w:t xml:space="preserve">  3    5</w:t><w:br/><w:t xml:space="preserve">  ----</w:t><w:br/><w:t>8|    |</w:t><w:br/><w:t>3|    |</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>')

$cell = $tbl.Cell(5, 1)
$p = $cell.Range.Paragraphs.Item(1).Range
$p.InsertXML('<?xml version="1.0" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:rPr><w:sz w:val="32"/></w:rPr><w:t>90 x 82</w:t><w:br/><w:t xml:space="preserve">  8    2</w:t><w:br/><w:t xml:space="preserve">  ----</w:t><w:br/><w:t>9|    |</w:t><w:br/><w:t>0|    |</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>')

$cell = $tbl.Cell(5, 2)
$p = $cell.Range.Paragraphs.Item(1).Range
$p.InsertXML('<?xml version="1.0" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:rPr><w:sz w:val="32"/></w:rPr><w:t>95 x 74</w:t><w:br/><w:t xml:space="preserve">  7    4</w:t><w:br/><w:t xml:space="preserve">  ----</w:t><w:br/><w:t>9|    |</w:t><w:br/><w:t>5|    |</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>')

$cell = $tbl.Cell(5, 3)
$p = $cell.Range.Paragraphs.Item(1).Range
$p.InsertXML('<?xml version="1.0" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:rPr><w:sz w:val="32"/></w:rPr><w:t>59 x 98</w:t><w:br/><w:t xml:space="preserve">  9    8</w:t><w:br/><w:t xml:space="preserve">  ----</w:t><w:br/><w:t>5|    |</w:t><w:br/><w:t>9|    |</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>')
